# Generate Report for Handoff
#
# The localization-status report tracks, per source file and per target
# language, whether a file has been handed off for translation yet. Two
# files (0aa50b45-...md is already in-flight, 8fc6fbd0-...md and
# 9a31e08b-...md were pending) just had their handoff generated, so:
#   - the "Latest Handoff Datetime" for 0aa50b45-...md / 8fc6fbd0-...md /
#     9a31e08b-...md is stamped with the real handoff time (per language
#     sheet), and
#   - the "Status" for 8fc6fbd0-...md moves from "In Translation" to
#     "Ready for handoff" (9a31e08b-...md was already "Ready for handoff").
#   - the Overview roll-up sheet reflects the same status change.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: roll-up status per language column -------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"

# --- zh-cn sheet: handoff datetime + status --------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D6").Value = "2016-03-09 06:28:12"
$wsZhCn.Range("B9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "2016-03-09 06:28:12"
$wsZhCn.Range("D10").Value = "2016-03-09 06:28:12"

# --- de-de sheet: handoff datetime + status --------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D6").Value = "2016-03-09 06:28:16"
$wsDeDe.Range("B9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "2016-03-09 06:28:16"
$wsDeDe.Range("D10").Value = "2016-03-09 06:28:16"
